# Auto-generated edit script: apply scheduled-runner market-price refresh
# to the Typhon_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2211.35
$ws.Range("J106").Value = 2911.0908
$ws.Range("L106").Value = 2911.0908
$ws.Range("N106").Value = -4173.0908
$ws.Range("H137").Value = 19050.965
$ws.Range("I137").Value = 1831.921
$ws.Range("K137").Value = 5495.763
$ws.Range("M137").Value = -2945.763

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1657.579
$ws.Range("I2").Value = 1030.0625
$ws.Range("J2").Value = 5004.3335
$ws.Range("K2").Value = 1030.0625
$ws.Range("L2").Value = 5004.3335
$ws.Range("M2").Value = -917.0625
$ws.Range("N2").Value = -5230.3335
$ws.Range("H5").Value = 1100
$ws.Range("I5").Value = 1100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -988
$ws.Range("N5").ClearContents()
$ws.Range("H61").Value = 2210.0967
$ws.Range("I61").Value = 1130.3125
$ws.Range("J61").Value = 3361.8667
$ws.Range("K61").Value = 1130.3125
$ws.Range("L61").Value = 3361.8667
$ws.Range("M61").Value = -918.3125
$ws.Range("N61").Value = -3785.8667
$ws.Range("H74").Value = 1152.52
$ws.Range("I74").Value = 848.7857
$ws.Range("J74").Value = 1539.091
$ws.Range("K74").Value = 848.7857
$ws.Range("L74").Value = 1539.091
$ws.Range("M74").Value = 25.21429999999998
$ws.Range("N74").Value = -3287.091
$ws.Range("H77").Value = 1152.52
$ws.Range("I77").Value = 848.7857
$ws.Range("J77").Value = 1539.091
$ws.Range("K77").Value = 4243.9285
$ws.Range("L77").Value = 7695.455
$ws.Range("M77").Value = 124.0715
$ws.Range("N77").Value = -16431.455
$ws.Range("H110").Value = 2381.7273
$ws.Range("I110").Value = 2119.9
$ws.Range("K110").Value = 2119.9
$ws.Range("M110").Value = -74.90000000000009
$ws.Range("H116").Value = 1657.579
$ws.Range("I116").Value = 1030.0625
$ws.Range("J116").Value = 5004.3335
$ws.Range("K116").Value = 1030.0625
$ws.Range("L116").Value = 5004.3335
$ws.Range("M116").Value = 1263.9375
$ws.Range("N116").Value = -9592.333500000001
$ws.Range("H136").Value = 2210.0967
$ws.Range("I136").Value = 1130.3125
$ws.Range("J136").Value = 3361.8667
$ws.Range("K136").Value = 3390.9375
$ws.Range("L136").Value = 10085.6001
$ws.Range("M136").Value = -840.9375
$ws.Range("N136").Value = -15185.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1657.579
$ws.Range("I3").Value = 1030.0625
$ws.Range("J3").Value = 5004.3335
$ws.Range("K3").Value = 1030.0625
$ws.Range("L3").Value = 5004.3335
$ws.Range("M3").Value = -916.0625
$ws.Range("N3").Value = -5232.3335
$ws.Range("H4").Value = 1100
$ws.Range("I4").Value = 1100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -985
$ws.Range("N4").ClearContents()
$ws.Range("H81").Value = 19270.285
$ws.Range("J81").Value = 19270.285
$ws.Range("L81").Value = 19270.285
$ws.Range("N81").Value = -21392.285
$ws.Range("H84").Value = 19270.285
$ws.Range("J84").Value = 19270.285
$ws.Range("L84").Value = 57810.855
$ws.Range("N84").Value = -68418.855
$ws.Range("H99").Value = 1999.8
$ws.Range("I99").Value = 1750
$ws.Range("K99").Value = 1750
$ws.Range("M99").Value = -252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13615.793
$ws.Range("J31").Value = 3950.3333
$ws.Range("L31").Value = 3950.3333
$ws.Range("N31").Value = -4540.3333
$ws.Range("H34").Value = 13615.793
$ws.Range("J34").Value = 3950.3333
$ws.Range("L34").Value = 3950.3333
$ws.Range("N34").Value = -4354.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5126.577
$ws.Range("J68").Value = 8393.532999999999
$ws.Range("L68").Value = 25180.599
$ws.Range("N68").Value = -26802.599
$ws.Range("H71").Value = 5126.577
$ws.Range("J71").Value = 8393.532999999999
$ws.Range("L71").Value = 75541.79699999999
$ws.Range("N71").Value = -83653.79699999999
$ws.Range("H107").Value = 4524.154
$ws.Range("I107").Value = 14780.857
$ws.Range("J107").Value = 745.3684
$ws.Range("K107").Value = 44342.571
$ws.Range("L107").Value = 2236.1052
$ws.Range("M107").Value = -42422.571
$ws.Range("N107").Value = -6076.1052
$ws.Range("H122").Value = 535.0833
$ws.Range("I122").Value = 387.5
$ws.Range("J122").Value = 608.875
$ws.Range("K122").Value = 3487.5
$ws.Range("L122").Value = 5479.875
$ws.Range("M122").Value = -1037.5
$ws.Range("N122").Value = -10379.875
$ws.Range("H131").Value = 107195.08
$ws.Range("J131").Value = 114452.47
$ws.Range("L131").Value = 343357.41
$ws.Range("N131").Value = -353437.41
$ws.Range("H132").Value = 926.2857
$ws.Range("I132").Value = 997.25
$ws.Range("J132").Value = 831.6667
$ws.Range("K132").Value = 8975.25
$ws.Range("L132").Value = 7485.0003
$ws.Range("M132").Value = -6445.25
$ws.Range("N132").Value = -12545.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2381.0667
$ws.Range("I97").Value = 1218.7273
$ws.Range("J97").Value = 5577.5
$ws.Range("K97").Value = 1218.7273
$ws.Range("L97").Value = 5577.5
$ws.Range("M97").Value = -722.7273
$ws.Range("N97").Value = -6569.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 326.6
$ws.Range("I9").Value = 280.5
$ws.Range("J9").Value = 511
$ws.Range("K9").Value = 280.5
$ws.Range("L9").Value = 511
$ws.Range("M9").Value = -56.5
$ws.Range("N9").Value = -959
$ws.Range("H22").Value = 2348.8333
$ws.Range("I22").Value = 5200.5
$ws.Range("J22").Value = 923
$ws.Range("K22").Value = 5200.5
$ws.Range("L22").Value = 923
$ws.Range("M22").Value = -4905.5
$ws.Range("N22").Value = -1513
$ws.Range("H27").Value = 2348.8333
$ws.Range("I27").Value = 5200.5
$ws.Range("J27").Value = 923
$ws.Range("K27").Value = 5200.5
$ws.Range("L27").Value = 923
$ws.Range("M27").Value = -5093.5
$ws.Range("N27").Value = -1137
$ws.Range("H100").Value = 2503
$ws.Range("I100").Value = 2503
$ws.Range("K100").Value = 2503
$ws.Range("M100").Value = -1962
$ws.Range("H136").Value = 2373.0605
$ws.Range("I136").Value = 1344.8334
$ws.Range("J136").Value = 3606.9333
$ws.Range("K136").Value = 4034.5002
$ws.Range("L136").Value = 10820.7999
$ws.Range("M136").Value = -1484.5002
$ws.Range("N136").Value = -15920.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1241.3125
$ws.Range("I136").Value = 860.1177
$ws.Range("J136").Value = 1673.3334
$ws.Range("K136").Value = 2580.3531
$ws.Range("L136").Value = 5020.0002
$ws.Range("M136").Value = -30.35310000000027
$ws.Range("N136").Value = -10120.0002

Write-Output "Applied 172 cell updates and 2 clears."
